$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CONFIG1 row (row 2): I2, J2, K2 bits change -> 0x91
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1

# Update CONFIG2 row (row 3): J3 bit change -> 0xD3
$ws.Range("J3").Value = 1

# Update CH3SET row (row 8): D8 bit change -> 0x65
$ws.Range("D8").Value = 0

# Update the view: scroll back to top-left and move selection to D9
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D9").Select() | Out-Null
